$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ======================================================================
# Sprint 21 upgrade: append new DB change-log rows (9-17) to the sheet.
# ======================================================================

# --- Row 9 ---
$ws.Cells.Item(9,1).Value2 = 40879
$ws.Cells.Item(9,2).Value2 = "AM"
$ws.Cells.Item(9,3).Value2 = "application.application_status_type"
$ws.Cells.Item(9,4).Value2 = "Change display_value for the dead status to Annulled"
$ws.Cells.Item(9,5).Value2 = "Dead is not a good description for a status. Annulled / Annulment accurately reflects the state of the application does not have the negative conotations of Dead. Agreed with Neil. "
$ws.Cells.Item(9,6).Value2 = 40891

# --- Row 10 ---
$ws.Cells.Item(10,1).Value2 = 40884
$ws.Cells.Item(10,2).Value2 = "Alexander"
$ws.Cells.Item(10,3).Value2 = "system.approle"
$ws.Cells.Item(10,4).Value2 = "Add new role - ""ManageBR"" with displayValue=""Manage business rules"" and Description=""Allows to manage business rules""."
$ws.Cells.Item(10,6).Value2 = 40891

# --- Row 11 ---
$ws.Cells.Item(11,1).Value2 = 40891
$ws.Cells.Item(11,2).Value2 = "AM"
$ws.Cells.Item(11,3).Value2 = "system.query_fields"
$ws.Cells.Item(11,4).Value2 = "Add field_name to query_fields table. Must be not null and match the name of the field used in the dynamic query."
$ws.Cells.Item(11,5).Value2 = "MyBatis does necessarily  not return the query results in the order they are selected and omits null values completely from the reuslt Map. Need to use the field_name to match the returned values to there field display values.  "
$ws.Cells.Item(11,6).Value2 = 40891
$ws.Cells.Item(11,7).Value2 = "Elton: I added column name instead of field_name, because I have used everywhere in db"

# --- Row 12 ---
$ws.Cells.Item(12,1).Value2 = 40891
$ws.Cells.Item(12,2).Value2 = "AM"
$ws.Cells.Item(12,3).Value2 = "system.query_fields data"
$ws.Cells.Item(12,4).Value2 = "Remove config data for the id field and the_geom field. Code now treats the values in query_fields as those requiring display names. All other fields in the query are obtained from the result set.  Ensure index is reordered from 0"
$ws.Cells.Item(12,5).Value2 = "The Mybatis code will merge the fields from query_fields with the fields in the query result to ensure all fields are included in the generic result. This minimises the configuration of query fields to only those that require localized display names. "
$ws.Cells.Item(12,6).Value2 = 40891

# --- Row 13 ---
$ws.Cells.Item(13,1).Value2 = 40891
$ws.Cells.Item(13,2).Value2 = "AM"
$ws.Cells.Item(13,3).Value2 = "system.query"
$ws.Cells.Item(13,4).Value2 = "Update all informationtool selects to use the_geom at the alias for st_asewkb"
$ws.Cells.Item(13,5).Value2 = "Ensures  the field name used for st_asewkb in the Generic Result  is the_geom. This avoids the need for additional config in query_fields for this field. "
$ws.Cells.Item(13,6).Value2 = 40891

# --- Row 14 ---
$ws.Cells.Item(14,1).Value2 = 40891
$ws.Cells.Item(14,2).Value2 = "AM"
$ws.Cells.Item(14,3).Value2 = "system.query"
$ws.Cells.Item(14,4).Value2 = "Update informationtool . Get_application query to select nr as label rather than label"
$ws.Cells.Item(14,5).Value2 = "label does not exist in the application table. "
$ws.Cells.Item(14,6).Value2 = 40891
$ws.Cells.Item(14,7).Value2 = "We can also leave it nr in this case without putting as label."

# --- Row 15 ---
$ws.Cells.Item(15,1).Value2 = 40891
$ws.Cells.Item(15,2).Value2 = "Elton"
$ws.Cells.Item(15,3).Value2 = "system.query_fields"
$ws.Cells.Item(15,4).Value2 = "Has to be changed to system.query_field "
$ws.Cells.Item(15,5).Value2 = "because everywhere in the db / model the naming is in singular"
$ws.Cells.Item(15,6).Value2 = 40891

# --- Row 16 ---
$ws.Cells.Item(16,1).Value2 = 40891
$ws.Cells.Item(16,2).Value2 = "Elton"
$ws.Cells.Item(16,3).Value2 = "Everywhere"
$ws.Cells.Item(16,4).Value2 = "triggers to handle default values for columns that are not optional but supplied with nulls has to be removed."
$ws.Cells.Item(16,5).Value2 = "This is handled in the service layer"
$ws.Cells.Item(16,6).Value2 = 40891

# --- Row 17 ---
$ws.Cells.Item(17,1).Value2 = 40891
$ws.Cells.Item(17,2).Value2 = "Alexander"
$ws.Cells.Item(17,3).Value2 = "party.is_rightholder"
$ws.Cells.Item(17,4).Value2 = "Return true if party connected with any RRR"
$ws.Cells.Item(17,5).Value2 = "Needed to distinguish parties as rightholders"
$ws.Cells.Item(17,6).Value2 = 40893
$ws.Cells.Item(17,7).Value2 = "CREATE OR REPLACE FUNCTION ""party"".is_rightholder(id character varying)
  RETURNS boolean AS
`$BODY`$
BEGIN
  return (SELECT (CASE (SELECT COUNT(1) FROM administrative.party_for_rrr ap WHERE ap.party_id = id) WHEN 0 THEN false ELSE true END));
END;
`$BODY`$
  LANGUAGE plpgsql VOLATILE
  COST 100;
ALTER FUNCTION ""party"".is_rightholder(character varying) OWNER TO postgres;"

# ----------------------------------------------------------------------
# Number formats: column A normally "d-mmm-yy" (style 3); rows 10 & 17
# use "mm-dd-yy" (=> built-in numFmtId 14) with left/top/wrap alignment.
# Column F uses "d-mmm" (=> built-in numFmtId 16) with top/wrap alignment
# for every new row (style 8).
# ----------------------------------------------------------------------

$ws.Cells.Item(9,1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(9,1).HorizontalAlignment = -4131
$ws.Cells.Item(9,1).VerticalAlignment = -4160
$ws.Cells.Item(9,1).WrapText = $true
$f = $ws.Cells.Item(9,6)
$f.NumberFormat = "d-mmm"
$f.VerticalAlignment = -4160
$f.WrapText = $true

$a = $ws.Cells.Item(10,1)
$a.NumberFormat = "mm-dd-yy"
$a.HorizontalAlignment = -4131
$a.VerticalAlignment = -4160
$a.WrapText = $true
$f = $ws.Cells.Item(10,6)
$f.NumberFormat = "d-mmm"
$f.VerticalAlignment = -4160
$f.WrapText = $true

$ws.Cells.Item(11,1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(11,1).HorizontalAlignment = -4131
$ws.Cells.Item(11,1).VerticalAlignment = -4160
$ws.Cells.Item(11,1).WrapText = $true
$f = $ws.Cells.Item(11,6)
$f.NumberFormat = "d-mmm"
$f.VerticalAlignment = -4160
$f.WrapText = $true

$ws.Cells.Item(12,1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(12,1).HorizontalAlignment = -4131
$ws.Cells.Item(12,1).VerticalAlignment = -4160
$ws.Cells.Item(12,1).WrapText = $true
$f = $ws.Cells.Item(12,6)
$f.NumberFormat = "d-mmm"
$f.VerticalAlignment = -4160
$f.WrapText = $true

$ws.Cells.Item(13,1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(13,1).HorizontalAlignment = -4131
$ws.Cells.Item(13,1).VerticalAlignment = -4160
$ws.Cells.Item(13,1).WrapText = $true
$f = $ws.Cells.Item(13,6)
$f.NumberFormat = "d-mmm"
$f.VerticalAlignment = -4160
$f.WrapText = $true

$ws.Cells.Item(14,1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(14,1).HorizontalAlignment = -4131
$ws.Cells.Item(14,1).VerticalAlignment = -4160
$ws.Cells.Item(14,1).WrapText = $true
$f = $ws.Cells.Item(14,6)
$f.NumberFormat = "d-mmm"
$f.VerticalAlignment = -4160
$f.WrapText = $true

$ws.Cells.Item(15,1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(15,1).HorizontalAlignment = -4131
$ws.Cells.Item(15,1).VerticalAlignment = -4160
$ws.Cells.Item(15,1).WrapText = $true
$f = $ws.Cells.Item(15,6)
$f.NumberFormat = "d-mmm"
$f.VerticalAlignment = -4160
$f.WrapText = $true

$ws.Cells.Item(16,1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(16,1).HorizontalAlignment = -4131
$ws.Cells.Item(16,1).VerticalAlignment = -4160
$ws.Cells.Item(16,1).WrapText = $true
$f = $ws.Cells.Item(16,6)
$f.NumberFormat = "d-mmm"
$f.VerticalAlignment = -4160
$f.WrapText = $true

$a = $ws.Cells.Item(17,1)
$a.NumberFormat = "mm-dd-yy"
$a.HorizontalAlignment = -4131
$a.VerticalAlignment = -4160
$a.WrapText = $true
$f = $ws.Cells.Item(17,6)
$f.NumberFormat = "d-mmm"
$f.VerticalAlignment = -4160
$f.WrapText = $true

# ----------------------------------------------------------------------
# Row heights (match the wrapped-text heights recorded by Excel).
# ----------------------------------------------------------------------

$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 105
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(14).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 255

# ----------------------------------------------------------------------
# Selection / view state.
# ----------------------------------------------------------------------
$ws.Range("G17").Select()
